$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price / volume / name / link data scraped on 2023-02-14
$updates = @{
    'D2' = '291.87'
    'D3' = '40.27'
    'E3' = '-0.44%'
    'D4' = '5.005'
    'E4' = '-0.89%'
    'D5' = '0.07289'
    'E5' = '-1.69%'
    'B6' = 'FTXToken'
    'C6' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D6' = '1.556'
    'E6' = '-1.97%'
    'B7' = 'MXToken'
    'C7' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D7' = '0.9258'
    'E7' = '-0.01%'
    'B8' = 'BTSEToken'
    'C8' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D8' = '2.361'
    'E8' = '-2.44%'
    'B9' = 'LiechtensteinCryptoassetsExchange'
    'C9' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D9' = '0.1158'
    'E9' = '-2.14%'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '0.1765'
    'E10' = '1.06%'
    'B11' = 'BitrueCoin'
    'C11' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D11' = '0.04361'
    'E11' = '3.75%'
    'B12' = 'MandalaExchangeToken'
    'C12' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D12' = '0.08744'
    'E12' = '0.04%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.1052'
    'E13' = '-0.24%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001274'
    'E14' = '-0.01%'
    'B15' = 'CoinExToken'
    'C15' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D15' = '0.03919'
    'E15' = '1.34%'
    'D16' = '0.006015'
    'E16' = '2.06%'
    'E17' = '-0.63%'
    'B18' = 'GateToken'
    'C18' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D18' = '4.282'
    'E18' = '-1.16%'
    'B19' = 'BitpandaEcosystemToken'
    'C19' = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
    'D19' = '0.3287'
    'E19' = '-1.86%'
    'B20' = 'MCDex'
    'C20' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D20' = '7.967'
    'E20' = '3.71%'
    'B21' = 'ProBitToken'
    'C21' = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
    'D21' = '0.1391'
    'E21' = '2.08%'
    'B22' = 'ZBToken'
    'C22' = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
    'D22' = '0.2774'
    'E22' = '-1.86%'
    'E23' = '-2.81%'
    'D24' = '0.003672'
    'E24' = '4.65%'
    'D25' = '0.0001201'
    'E25' = '-8.30%'
    'D26' = '0.0003726'
    'E26' = '-0.83%'
    'E38' = '0.13%'
    'D39' = '0.05073'
    'E39' = '1.38%'
    'D40' = '0.005665'
    'E40' = '35.89%'
    'D41' = '0.007855'
    'E41' = '0.82%'
    'E42' = '0.71%'
    'D43' = '0.007392'
    'E43' = '-0.39%'
    'D44' = '0.007275'
    'E44' = '1.68%'
    'D45' = '0.3191'
    'E45' = '0.20%'
    'D46' = '0.00006187'
    'E46' = '-7.96%'
    'E47' = '-0.82%'
    'D48' = '0.04955'
    'E48' = '-80.32%'
    'D49' = '0.00002102'
    'E49' = '-0.82%'
    'D50' = '0.0002002'
    'E50' = '-0.82%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text interpretation so numeric-looking / percent-looking strings
    # (e.g. "291.87", "-0.44%") are stored verbatim as text, not coerced to numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Reset the style back to the workbook default so no stray formatting/style
    # index is introduced by the temporary text number-format.
    $cell.Style = "Normal"
}
